# Apply 2023-12-11 data update to violent-crime-full-year workbook
# Updates column J (year 2023) values across Citywide Totals, By Neighborhood,
# and each individual neighborhood sheet, per the published diff.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet='Citywide Totals'; Row=2; Value=7263}
    @{Sheet='Citywide Totals'; Row=3; Value=7648}
    @{Sheet='Citywide Totals'; Row=4; Value=1668}
    @{Sheet='Citywide Totals'; Row=5; Value=598}
    @{Sheet='Citywide Totals'; Row=6; Value=10415}
    @{Sheet='Citywide Totals'; Row=7; Value=27592}
    @{Sheet='By Neighborhood'; Row=2; Value=218}
    @{Sheet='By Neighborhood'; Row=4; Value=127}
    @{Sheet='By Neighborhood'; Row=5; Value=83}
    @{Sheet='By Neighborhood'; Row=7; Value=787}
    @{Sheet='By Neighborhood'; Row=8; Value=1739}
    @{Sheet='By Neighborhood'; Row=11; Value=496}
    @{Sheet='By Neighborhood'; Row=15; Value=342}
    @{Sheet='By Neighborhood'; Row=16; Value=110}
    @{Sheet='By Neighborhood'; Row=19; Value=793}
    @{Sheet='By Neighborhood'; Row=20; Value=590}
    @{Sheet='By Neighborhood'; Row=21; Value=78}
    @{Sheet='By Neighborhood'; Row=25; Value=141}
    @{Sheet='By Neighborhood'; Row=29; Value=1471}
    @{Sheet='By Neighborhood'; Row=30; Value=97}
    @{Sheet='By Neighborhood'; Row=31; Value=287}
    @{Sheet='By Neighborhood'; Row=33; Value=1251}
    @{Sheet='By Neighborhood'; Row=36; Value=371}
    @{Sheet='By Neighborhood'; Row=37; Value=843}
    @{Sheet='By Neighborhood'; Row=41; Value=208}
    @{Sheet='By Neighborhood'; Row=42; Value=1183}
    @{Sheet='By Neighborhood'; Row=43; Value=232}
    @{Sheet='By Neighborhood'; Row=44; Value=218}
    @{Sheet='By Neighborhood'; Row=48; Value=308}
    @{Sheet='By Neighborhood'; Row=52; Value=705}
    @{Sheet='By Neighborhood'; Row=53; Value=413}
    @{Sheet='By Neighborhood'; Row=54; Value=547}
    @{Sheet='By Neighborhood'; Row=55; Value=437}
    @{Sheet='By Neighborhood'; Row=57; Value=130}
    @{Sheet='By Neighborhood'; Row=59; Value=31}
    @{Sheet='By Neighborhood'; Row=60; Value=161}
    @{Sheet='By Neighborhood'; Row=63; Value=83}
    @{Sheet='By Neighborhood'; Row=65; Value=693}
    @{Sheet='By Neighborhood'; Row=67; Value=1022}
    @{Sheet='By Neighborhood'; Row=70; Value=40}
    @{Sheet='By Neighborhood'; Row=73; Value=267}
    @{Sheet='By Neighborhood'; Row=76; Value=391}
    @{Sheet='By Neighborhood'; Row=77; Value=193}
    @{Sheet='By Neighborhood'; Row=78; Value=320}
    @{Sheet='By Neighborhood'; Row=79; Value=757}
    @{Sheet='By Neighborhood'; Row=83; Value=551}
    @{Sheet='By Neighborhood'; Row=84; Value=228}
    @{Sheet='By Neighborhood'; Row=85; Value=1129}
    @{Sheet='By Neighborhood'; Row=88; Value=292}
    @{Sheet='By Neighborhood'; Row=89; Value=342}
    @{Sheet='By Neighborhood'; Row=90; Value=290}
    @{Sheet='By Neighborhood'; Row=93; Value=119}
    @{Sheet='By Neighborhood'; Row=94; Value=303}
    @{Sheet='By Neighborhood'; Row=95; Value=393}
    @{Sheet='By Neighborhood'; Row=96; Value=308}
    @{Sheet='By Neighborhood'; Row=98; Value=204}
    @{Sheet='By Neighborhood'; Row=99; Value=423}
    @{Sheet='By Neighborhood'; Row=101; Value=27592}
    @{Sheet='West Ridge'; Row=6; Value=115}
    @{Sheet='West Ridge'; Row=7; Value=308}
    @{Sheet='Auburn Gresham'; Row=2; Value=249}
    @{Sheet='Auburn Gresham'; Row=6; Value=250}
    @{Sheet='Auburn Gresham'; Row=7; Value=787}
    @{Sheet='Belmont Cragin'; Row=3; Value=85}
    @{Sheet='Belmont Cragin'; Row=6; Value=235}
    @{Sheet='Belmont Cragin'; Row=7; Value=496}
    @{Sheet='Uptown'; Row=4; Value=34}
    @{Sheet='Uptown'; Row=6; Value=105}
    @{Sheet='Uptown'; Row=7; Value=342}
    @{Sheet='South Shore'; Row=2; Value=303}
    @{Sheet='South Shore'; Row=6; Value=321}
    @{Sheet='South Shore'; Row=7; Value=1129}
    @{Sheet='Little Village'; Row=2; Value=165}
    @{Sheet='Little Village'; Row=6; Value=304}
    @{Sheet='Little Village'; Row=7; Value=705}
    @{Sheet='Logan Square'; Row=2; Value=74}
    @{Sheet='Logan Square'; Row=3; Value=51}
    @{Sheet='Logan Square'; Row=7; Value=413}
    @{Sheet='Austin'; Row=2; Value=458}
    @{Sheet='Austin'; Row=4; Value=92}
    @{Sheet='Austin'; Row=6; Value=640}
    @{Sheet='Austin'; Row=7; Value=1739}
    @{Sheet='South Chicago'; Row=2; Value=164}
    @{Sheet='South Chicago'; Row=7; Value=551}
    @{Sheet='Garfield Park'; Row=2; Value=282}
    @{Sheet='Garfield Park'; Row=5; Value=53}
    @{Sheet='Garfield Park'; Row=6; Value=446}
    @{Sheet='Garfield Park'; Row=7; Value=1251}
    @{Sheet='West Pullman'; Row=3; Value=142}
    @{Sheet='West Pullman'; Row=7; Value=393}
    @{Sheet='Grand Crossing'; Row=6; Value=244}
    @{Sheet='Grand Crossing'; Row=7; Value=843}
    @{Sheet='New City'; Row=6; Value=258}
    @{Sheet='New City'; Row=7; Value=693}
    @{Sheet='Woodlawn'; Row=2; Value=116}
    @{Sheet='Woodlawn'; Row=7; Value=423}
    @{Sheet='Fuller Park'; Row=3; Value=38}
    @{Sheet='Fuller Park'; Row=7; Value=97}
    @{Sheet='Gage Park'; Row=3; Value=69}
    @{Sheet='Gage Park'; Row=6; Value=103}
    @{Sheet='Gage Park'; Row=7; Value=287}
    @{Sheet='North Lawndale'; Row=6; Value=282}
    @{Sheet='North Lawndale'; Row=7; Value=1022}
    @{Sheet='South Deering'; Row=2; Value=67}
    @{Sheet='South Deering'; Row=3; Value=73}
    @{Sheet='South Deering'; Row=7; Value=228}
    @{Sheet='Loop'; Row=3; Value=110}
    @{Sheet='Loop'; Row=6; Value=252}
    @{Sheet='Loop'; Row=7; Value=547}
    @{Sheet='Englewood'; Row=3; Value=520}
    @{Sheet='Englewood'; Row=6; Value=371}
    @{Sheet='Englewood'; Row=7; Value=1471}
    @{Sheet='Lake View'; Row=2; Value=52}
    @{Sheet='Lake View'; Row=7; Value=308}
    @{Sheet='Chatham'; Row=2; Value=194}
    @{Sheet='Chatham'; Row=3; Value=227}
    @{Sheet='Chatham'; Row=7; Value=793}
    @{Sheet='Irving Park'; Row=3; Value=50}
    @{Sheet='Irving Park'; Row=6; Value=87}
    @{Sheet='Irving Park'; Row=7; Value=218}
    @{Sheet='River North'; Row=3; Value=84}
    @{Sheet='River North'; Row=4; Value=30}
    @{Sheet='River North'; Row=7; Value=391}
    @{Sheet='Hermosa'; Row=3; Value=30}
    @{Sheet='Hermosa'; Row=7; Value=208}
    @{Sheet='Humboldt Park'; Row=2; Value=245}
    @{Sheet='Humboldt Park'; Row=3; Value=239}
    @{Sheet='Humboldt Park'; Row=7; Value=1183}
    @{Sheet='Rogers Park'; Row=2; Value=85}
    @{Sheet='Rogers Park'; Row=7; Value=320}
    @{Sheet='Lower West Side'; Row=2; Value=83}
    @{Sheet='Lower West Side'; Row=7; Value=437}
    @{Sheet='Chinatown'; Row=2; Value=13}
    @{Sheet='Chinatown'; Row=6; Value=52}
    @{Sheet='Chinatown'; Row=7; Value=78}
    @{Sheet='Roseland'; Row=6; Value=227}
    @{Sheet='Roseland'; Row=7; Value=757}
    @{Sheet='Chicago Lawn'; Row=2; Value=164}
    @{Sheet='Chicago Lawn'; Row=3; Value=196}
    @{Sheet='Chicago Lawn'; Row=5; Value=14}
    @{Sheet='Chicago Lawn'; Row=6; Value=172}
    @{Sheet='Chicago Lawn'; Row=7; Value=590}
    @{Sheet='Grand Boulevard'; Row=2; Value=120}
    @{Sheet='Grand Boulevard'; Row=3; Value=120}
    @{Sheet='Grand Boulevard'; Row=6; Value=111}
    @{Sheet='Grand Boulevard'; Row=7; Value=371}
    @{Sheet='West Lawn'; Row=3; Value=37}
    @{Sheet='West Lawn'; Row=7; Value=119}
    @{Sheet='West Loop'; Row=6; Value=159}
    @{Sheet='West Loop'; Row=7; Value=303}
    @{Sheet='East Side'; Row=3; Value=42}
    @{Sheet='East Side'; Row=7; Value=141}
    @{Sheet='Brighton Park'; Row=3; Value=69}
    @{Sheet='Brighton Park'; Row=7; Value=342}
    @{Sheet='Wicker Park'; Row=3; Value=29}
    @{Sheet='Wicker Park'; Row=7; Value=204}
    @{Sheet='Portage Park'; Row=2; Value=85}
    @{Sheet='Portage Park'; Row=6; Value=99}
    @{Sheet='Portage Park'; Row=7; Value=267}
    @{Sheet='Montclare'; Row=6; Value=7}
    @{Sheet='Montclare'; Row=7; Value=31}
    @{Sheet='Albany Park'; Row=6; Value=79}
    @{Sheet='Albany Park'; Row=7; Value=218}
    @{Sheet='O''Hare'; Row=6; Value=6}
    @{Sheet='O''Hare'; Row=7; Value=40}
    @{Sheet='United Center'; Row=2; Value=59}
    @{Sheet='United Center'; Row=7; Value=292}
    @{Sheet='Armour Square'; Row=3; Value=15}
    @{Sheet='Armour Square'; Row=6; Value=41}
    @{Sheet='Armour Square'; Row=7; Value=83}
    @{Sheet='Washington Heights'; Row=6; Value=87}
    @{Sheet='Washington Heights'; Row=7; Value=290}
    @{Sheet='Mckinley Park'; Row=2; Value=34}
    @{Sheet='Mckinley Park'; Row=7; Value=130}
    @{Sheet='Morgan Park'; Row=6; Value=49}
    @{Sheet='Morgan Park'; Row=7; Value=161}
    @{Sheet='Hyde Park'; Row=6; Value=137}
    @{Sheet='Hyde Park'; Row=7; Value=232}
    @{Sheet='Riverdale'; Row=2; Value=73}
    @{Sheet='Riverdale'; Row=7; Value=193}
    @{Sheet='Archer Heights'; Row=2; Value=38}
    @{Sheet='Archer Heights'; Row=7; Value=127}
    @{Sheet='Bucktown'; Row=2; Value=13}
    @{Sheet='Bucktown'; Row=7; Value=110}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Cells.Item($u.Row, 10).Value = $u.Value
}

Write-Host "Applied" $updates.Count "2023 (column J) updates for 2023-12-11 data refresh"
